$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 9 and 10 (pushes old rows 9-13 down to 11-15),
# inheriting formatting from the row that was previously at 9.
$ws.Rows("9:10").Insert()

# Row heights for the two new rows
$ws.Rows(9).RowHeight = 30
$ws.Rows(10).RowHeight = 30

# Row 9: PHU_Num
$ws.Range("A9").Value = "PHU_Num"
$ws.Range("B9").Value = "PHU_Num"
$ws.Range("C9").Value = "ID of the Public Health Unit (PHU) region in which the long term care home is located."
$ws.Range("D9").Value = "ID de la région du bureau de santé publique dans laquelle se trouve le foyer de soins de longue durée."

# Row 10: PHU
$ws.Range("A10").Value = "PHU"
$ws.Range("B10").Value = "PHU"
$ws.Range("C10").Value = "Name of the Public Health Unit (PHU) region in which the long term care home is located."
$ws.Range("D10").Value = "Nom de la région du bureau de santé publique dans laquelle se trouve le foyer de soins de longue durée."

# Formatting for C9:C10 and D9:D10 to match new cellXfs (vertical top + wrap text)
$ws.Range("C9:C10").WrapText = $true
$ws.Range("C9:C10").VerticalAlignment = -4160

$ws.Range("D9:D10").WrapText = $true
$ws.Range("D9:D10").VerticalAlignment = -4160
$ws.Range("D9:D10").Font.Color = 0

# Update the sheet view: scrolled down a bit, selection moved to D8
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D8").Select()
